$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "last f3"
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = "first f3"
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = "HCM"
$ws.Range("G3").Value = "add up f3"
